$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 527503
$ws.Range("I2").Value = 550.8889
$ws.Range("J2").Value = 1001759.9
$ws.Range("K2").Value = 550.8889
$ws.Range("L2").Value = 1001759.9
$ws.Range("M2").Value = -437.8889
$ws.Range("N2").Value = -1001985.9
$ws.Range("H40").Value = 50003084
$ws.Range("I40").Value = 3379.6
$ws.Range("J40").Value = 100002780
$ws.Range("K40").Value = 3379.6
$ws.Range("L40").Value = 100002780
$ws.Range("M40").Value = -3204.6
$ws.Range("N40").Value = -100003130
$ws.Range("H86").Value = 9077.637000000001
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 9077.637000000001
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 9077.637000000001
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -11323.637
$ws.Range("H89").Value = 9077.637000000001
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 9077.637000000001
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 45388.185
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -56620.185
$ws.Range("H99").Value = 3613.2222
$ws.Range("J99").Value = 4599.4287
$ws.Range("L99").Value = 13798.2861
$ws.Range("N99").Value = -16794.2861
$ws.Range("H100").Value = 12444
$ws.Range("I100").Value = 7799.3335
$ws.Range("K100").Value = 7799.3335
$ws.Range("M100").Value = -7258.3335
$ws.Range("H132").Value = 5809.2144
$ws.Range("I132").Value = 3277.0527
$ws.Range("K132").Value = 9831.158100000001
$ws.Range("M132").Value = -7301.158100000001
$ws.Range("H135").Value = 3793.6875
$ws.Range("J135").Value = 6779
$ws.Range("L135").Value = 61011
$ws.Range("N135").Value = -66081
$ws.Range("H137").Value = 2929.3333
$ws.Range("I137").Value = 2909.2856
$ws.Range("J137").Value = 2999.5
$ws.Range("K137").Value = 8727.856800000001
$ws.Range("L137").Value = 8998.5
$ws.Range("M137").Value = -6177.856800000001
$ws.Range("N137").Value = -14098.5
$ws.Range("H138").Value = 5674.6294
$ws.Range("I138").Value = 2997.55
$ws.Range("K138").Value = 8992.650000000001
$ws.Range("M138").Value = -3852.650000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 79989
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1655
$ws.Range("I94").Value = 1641.3125
$ws.Range("J94").Value = 1691.5
$ws.Range("K94").Value = 1641.3125
$ws.Range("L94").Value = 1691.5
$ws.Range("M94").Value = -1190.3125
$ws.Range("N94").Value = -2593.5
$ws.Range("H107").Value = 7416.5713
$ws.Range("I107").Value = 9202
$ws.Range("J107").Value = 2953
$ws.Range("K107").Value = 9202
$ws.Range("L107").Value = 2953
$ws.Range("M107").Value = -7282
$ws.Range("N107").Value = -6793

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29414902
$ws.Range("I31").Value = 52634750
$ws.Range("J31").Value = 3090.6667
$ws.Range("K31").Value = 52634750
$ws.Range("L31").Value = 3090.6667
$ws.Range("M31").Value = -52634455
$ws.Range("N31").Value = -3680.6667
$ws.Range("H34").Value = 29414902
$ws.Range("I34").Value = 52634750
$ws.Range("J34").Value = 3090.6667
$ws.Range("K34").Value = 52634750
$ws.Range("L34").Value = 3090.6667
$ws.Range("M34").Value = -52634548
$ws.Range("N34").Value = -3494.6667
$ws.Range("H122").Value = 5430.1113
$ws.Range("I122").Value = 5782.2
$ws.Range("J122").Value = 4990
$ws.Range("K122").Value = 17346.6
$ws.Range("L122").Value = 14970
$ws.Range("M122").Value = -14896.6
$ws.Range("N122").Value = -19870

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1994
$ws.Range("I5").Value = 1994
$ws.Range("K5").Value = 5982
$ws.Range("M5").Value = -5870
$ws.Range("H11").Value = 5003.8887
$ws.Range("I11").Value = 1472.4286
$ws.Range("J11").Value = 17364
$ws.Range("K11").Value = 4417.2858
$ws.Range("L11").Value = 52092
$ws.Range("M11").Value = -4277.2858
$ws.Range("N11").Value = -52372
$ws.Range("H14").Value = 8142.2856
$ws.Range("I14").Value = 8142.2856
$ws.Range("K14").Value = 24426.8568
$ws.Range("M14").Value = -24253.8568
$ws.Range("H128").Value = 250390
$ws.Range("I128").Value = 250390
$ws.Range("K128").Value = 751170
$ws.Range("M128").Value = -746190
$ws.Range("H135").Value = 1994
$ws.Range("I135").Value = 1994
$ws.Range("K135").Value = 17946
$ws.Range("M135").Value = -15411

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 554.43475
$ws.Range("I97").Value = 614.7778
$ws.Range("K97").Value = 614.7778
$ws.Range("M97").Value = -118.7778
$ws.Range("H102").Value = 3031.6785
$ws.Range("I102").Value = 3007.7827
$ws.Range("K102").Value = 3007.7827
$ws.Range("M102").Value = -1385.7827

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1326.2941
$ws.Range("I55").Value = 1302.7142
$ws.Range("J55").Value = 1342.8
$ws.Range("K55").Value = 1302.7142
$ws.Range("L55").Value = 1342.8
$ws.Range("M55").Value = -1129.7142
$ws.Range("N55").Value = -1688.8
$ws.Range("H68").Value = 4169486
$ws.Range("J68").Value = 3466.6667
$ws.Range("L68").Value = 3466.6667
$ws.Range("N68").Value = -4964.6667
$ws.Range("H71").Value = 4169486
$ws.Range("J71").Value = 3466.6667
$ws.Range("L71").Value = 17333.3335
$ws.Range("N71").Value = -24821.3335
$ws.Range("H82").Value = 3064.5
$ws.Range("I82").Value = 3359.4285
$ws.Range("J82").Value = 1000
$ws.Range("K82").Value = 3359.4285
$ws.Range("L82").Value = 1000
$ws.Range("M82").Value = -2998.4285
$ws.Range("N82").Value = -1722
$ws.Range("H85").Value = 3064.5
$ws.Range("I85").Value = 3359.4285
$ws.Range("J85").Value = 1000
$ws.Range("K85").Value = 3359.4285
$ws.Range("L85").Value = 1000
$ws.Range("M85").Value = -2111.4285
$ws.Range("N85").Value = -3496
$ws.Range("H93").Value = 2419559
$ws.Range("I93").Value = 2017.6111
$ws.Range("J93").Value = 11122708
$ws.Range("K93").Value = 2017.6111
$ws.Range("L93").Value = 11122708
$ws.Range("M93").Value = -769.6111000000001
$ws.Range("N93").Value = -11125204
$ws.Range("H132").Value = 4039.348
$ws.Range("I132").Value = 2693.7144
$ws.Range("K132").Value = 8081.1432
$ws.Range("M132").Value = -5551.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 57500
$ws.Range("J56").Value = 57500
$ws.Range("L56").Value = 57500
$ws.Range("N56").Value = -58928
$ws.Range("H96").Value = 10519.111
$ws.Range("I96").Value = 9096
$ws.Range("J96").Value = 15500
$ws.Range("K96").Value = 9096
$ws.Range("L96").Value = 15500
$ws.Range("M96").Value = -7723
$ws.Range("N96").Value = -18246
$ws.Range("H132").Value = 260297.03
$ws.Range("I132").Value = 3391.8965
$ws.Range("K132").Value = 10175.6895
$ws.Range("M132").Value = -7645.6895
$ws.Range("H136").Value = 323717.66
$ws.Range("I136").Value = 12002.654
$ws.Range("K136").Value = 36007.962
$ws.Range("M136").Value = -33457.962
